$d = $word.ActiveDocument

# The change targets the weekly-activities table (the 2nd top-level table
# in the document), in the row labelled "Dienstag", in the content column
# (column 2). That cell currently holds a single empty paragraph; we add
# text describing the day's activities: "Nutzwertanalyse, Kaufvertrag,
# Hardware" -- authored as three runs (matching how the text was typed /
# appended in separate edits), all using the same 11pt (half-point 22)
# font size already used throughout this table.
$table = $d.Tables.Item(2)
$cell = $table.Cell(3, 2)

$cellRange = $cell.Range
# Exclude the end-of-cell marker so we only touch the paragraph's own
# content (there is currently no text/run inside this paragraph).
$cellRange.End = $cellRange.End - 1
$startPos = $cellRange.Start

# --- Run 1: "Nutzwertanalyse, " ---
$cellRange.Text = "Nutzwertanalyse, "
$pos1 = $cellRange.End
$run1 = $d.Range($startPos, $pos1)
$run1.Font.Size = 11

# --- Run 2: "Kaufvertrag" ---
$run1.InsertAfter("Kaufvertrag")
$pos2 = $run1.End
$run2 = $d.Range($pos1, $pos2)
$run2.Font.Size = 11

# --- Run 3: ", Hardware" ---
$run2.InsertAfter(", Hardware")
$pos3 = $run2.End
$run3 = $d.Range($pos2, $pos3)
$run3.Font.Size = 11
